# Applies the diff: wraps the literal "burrow"/"pot" words that appear
# inside certain phrases with colored <env>/<tl> XML-tag markup runs
# (Courier New, blue, size 9) while keeping the tag contents (the actual
# word) in plain black formatting, matching the style already used
# elsewhere in the document for similar tags.

$d = $word.ActiveDocument

function Set-TagFont($range) {
    $range.Font.Name = "Courier New"
    $range.Font.Color = 16711680   # RGB(0,0,255) == w:val 0000ff in COM BGR encoding
    $range.Font.Size = 9           # sz/szCs = 18 (half-points) == 9pt
}

function Set-PlainFont($range) {
    $range.Font.Name = "Arial"
    $range.Font.Color = 0          # RGB(0,0,0) == w:val 000000
}

# --- Edit 1: "For making <al>rabbits</al> come out of a burrow" -------
# " come out of a burrow"  ->  " come out of a " + <env> + "burrow" + </env>
$rng = $d.Content
$rng.Find.Execute("come out of a burrow", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
$end = $rng.End
$pre = "come out of a "
$wordStart = $start + $pre.Length

$rWord = $d.Range($wordStart, $end)
$rWord.InsertBefore("<env>")
# recompute: the "burrow" word range now sits right after the inserted "<env>"
$tagOpenStart = $wordStart
$tagOpenEnd = $tagOpenStart + "<env>".Length
$wordStart2 = $tagOpenEnd
$wordEnd2 = $wordStart2 + "burrow".Length

Set-TagFont ($d.Range($tagOpenStart, $tagOpenEnd))

$rAfterWord = $d.Range($wordEnd2, $wordEnd2)
$rAfterWord.InsertBefore("</env>")
$tagCloseStart = $wordEnd2
$tagCloseEnd = $tagCloseStart + "</env>".Length
Set-TagFont ($d.Range($tagCloseStart, $tagCloseEnd))

# --- Edit 2: "Take some <m>embers</m> in a pot, &amp; having put" -----
# " in a pot, "  ->  " in a " + <tl> + "pot" + </tl> + ", "
$rng2 = $d.Content
$rng2.Find.Execute("in a pot, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start
$end2 = $rng2.End
$pre2 = "in a "
$wStart = $start2 + $pre2.Length
$wEnd = $wStart + "pot".Length

$rPot = $d.Range($wStart, $wEnd)
$rPot.InsertBefore("<tl>")
$openStart2 = $wStart
$openEnd2 = $openStart2 + "<tl>".Length
$potStart2 = $openEnd2
$potEnd2 = $potStart2 + "pot".Length
Set-TagFont ($d.Range($openStart2, $openEnd2))

$rAfterPot = $d.Range($potEnd2, $potEnd2)
$rAfterPot.InsertBefore("</tl>")
$closeStart2 = $potEnd2
$closeEnd2 = $closeStart2 + "</tl>".Length
Set-TagFont ($d.Range($closeStart2, $closeEnd2))

# --- Edit 3: "put it in the burrows &amp; block it" -------------------
# " the burrows"  ->  " the " + <env> + "burrow" + "s" + </env>
$rng3 = $d.Content
$rng3.Find.Execute("in the burrows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start3 = $rng3.Start
$end3 = $rng3.End
$pre3 = "in the "
$wStart3 = $start3 + $pre3.Length
$wEnd3 = $wStart3 + "burrow".Length  # the existing "s" run remains untouched after this

$rBurrow3 = $d.Range($wStart3, $wEnd3)
$rBurrow3.InsertBefore("<env>")
$openStart3 = $wStart3
$openEnd3 = $openStart3 + "<env>".Length
$burrowStart3 = $openEnd3
$burrowEnd3 = $burrowStart3 + "burrow".Length
Set-TagFont ($d.Range($openStart3, $openEnd3))

# the "s" sits right after burrowEnd3 (pre-existing run, untouched),
# so the closing tag must be inserted after that "s" character.
$sEnd3 = $burrowEnd3 + "s".Length
$rAfterS = $d.Range($sEnd3, $sEnd3)
$rAfterS.InsertBefore("</env>")
$closeStart3 = $sEnd3
$closeEnd3 = $closeStart3 + "</env>".Length
Set-TagFont ($d.Range($closeStart3, $closeEnd3))
